$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update column F (想去人数 / "want to go" count)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 471
$wsExhibit.Range("F3").Value = 5655
$wsExhibit.Range("F4").Value = 389
$wsExhibit.Range("F6").Value = 89

# Sheet "全部类型" (All Types) - same underlying rows, update column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 471
$wsAll.Range("F3").Value = 5655
$wsAll.Range("F4").Value = 389
$wsAll.Range("F7").Value = 89
